$wb = $excel.ActiveWorkbook

# Sheet index 1 (sheet1, "展览")
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F8").Value = 2262
$ws1.Range("F9").Value = 1484
$ws1.Range("F11").Value = 628
$ws1.Range("F13").Value = 2596
$ws1.Range("F15").Value = 1431
$ws1.Range("F16").Value = 5671
$ws1.Range("F17").Value = 9
$ws1.Range("F18").Value = 5463
$ws1.Range("F19").Value = 2058
$ws1.Range("F20").Value = 2961
$ws1.Range("F21").Value = 3394
$ws1.Range("F22").Value = 198
$ws1.Range("F23").Value = 1662
$ws1.Range("F24").Value = 34
$ws1.Range("F25").Value = 279
$ws1.Range("F26").Value = 854
$ws1.Range("F27").Value = 152
$ws1.Range("F28").Value = 13
$ws1.Range("F30").Value = 1060
$ws1.Range("F31").Value = 2237
$ws1.Range("F33").Value = 136
$ws1.Range("F34").Value = 317
$ws1.Range("F35").Value = 839
$ws1.Range("F37").Value = 403
$ws1.Range("F38").Value = 468

# Sheet index 2 (sheet2, "演出")
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F12").Value = 28

# Sheet index 4 (sheet4, "全部类型")
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F12").Value = 2262
$ws4.Range("F13").Value = 1484
$ws4.Range("F15").Value = 629
$ws4.Range("F17").Value = 28
$ws4.Range("F18").Value = 2596
$ws4.Range("F19").Value = 1431
$ws4.Range("F24").Value = 5671
$ws4.Range("F25").Value = 9
$ws4.Range("F26").Value = 5463
$ws4.Range("F27").Value = 2058
$ws4.Range("F28").Value = 2961
$ws4.Range("F29").Value = 3394
$ws4.Range("F31").Value = 198
$ws4.Range("F34").Value = 1662
$ws4.Range("F36").Value = 279
$ws4.Range("F37").Value = 854
$ws4.Range("F38").Value = 152
$ws4.Range("F39").Value = 13
$ws4.Range("F42").Value = 2237
$ws4.Range("F44").Value = 136
$ws4.Range("F45").Value = 317
$ws4.Range("F46").Value = 839
$ws4.Range("F48").Value = 403
$ws4.Range("F49").Value = 468
